$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data rows 2-6 with corrected IFRS figures ---
# Row 2
$ws.Range("D2").Value = 169
$ws.Range("E2").Value = 19
$ws.Range("F2").Value = 18
$ws.Range("G2").Value = 17
$ws.Range("H2").Value = 19
$ws.Range("I2").Value = 20
$ws.Range("J2").Value = -2
$ws.Range("K2").Value = 440
$ws.Range("L2").Value = 197
$ws.Range("M2").Value = 244
$ws.Range("N2").Value = 245
$ws.Range("O2").Value = -1
$ws.Range("P2").Value = 172
$ws.Range("Q2").Value = -21
$ws.Range("R2").Value = 14
$ws.Range("S2").Value = 49
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = -21
$ws.Range("V2").Value = 117
$ws.Range("W2").Value = 11.05
$ws.Range("X2").Value = 11.06
$ws.Range("Y2").Value = 8.699999999999999
$ws.Range("Z2").Value = 4.56
$ws.Range("AA2").Value = 80.7
$ws.Range("AB2").Value = 61.54
$ws.Range("AC2").Value = 59
$ws.Range("AD2").Value = 20.57
$ws.Range("AE2").Value = 846
$ws.Range("AF2").Value = 1.44
$ws.Range("AG2").Value = 15
$ws.Range("AH2").Value = 1.23
$ws.Range("AI2").Value = 21.32
$ws.Range("AJ2").Value = 34458589

# Row 3
$ws.Range("D3").Value = 172
$ws.Range("E3").Value = -34
$ws.Range("F3").Value = -34
$ws.Range("G3").Value = -30
$ws.Range("H3").Value = -6
$ws.Range("I3").Value = 8
$ws.Range("J3").Value = -15
$ws.Range("K3").Value = 676
$ws.Range("L3").Value = 336
$ws.Range("M3").Value = 340
$ws.Range("N3").Value = 354
$ws.Range("O3").Value = -15
$ws.Range("P3").Value = 172
$ws.Range("Q3").Value = -108
$ws.Range("R3").Value = -25
$ws.Range("S3").Value = 269
$ws.Range("T3").Value = 3
$ws.Range("U3").Value = -111
$ws.Range("V3").Value = 274
$ws.Range("W3").Value = -19.51
$ws.Range("X3").Value = -3.7
$ws.Range("Y3").Value = 2.74
$ws.Range("Z3").Value = -1.14
$ws.Range("AA3").Value = 98.95999999999999
$ws.Range("AB3").Value = 107.68
$ws.Range("AC3").Value = 24
$ws.Range("AD3").Value = 216.75
$ws.Range("AE3").Value = 1083
$ws.Range("AF3").Value = 4.77
$ws.Range("AG3").Value = 16
$ws.Range("AH3").Value = 0.31
$ws.Range("AI3").Value = 63.84
$ws.Range("AJ3").Value = 34458589

# Row 4
$ws.Range("D4").Value = 505
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 8
$ws.Range("G4").Value = 80
$ws.Range("H4").Value = 43
$ws.Range("I4").Value = 42
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 823
$ws.Range("L4").Value = 436
$ws.Range("M4").Value = 387
$ws.Range("N4").Value = 396
$ws.Range("O4").Value = -9
$ws.Range("P4").Value = 172
$ws.Range("Q4").Value = 111
$ws.Range("R4").Value = -29
$ws.Range("S4").Value = -128
$ws.Range("T4").Value = 3
$ws.Range("U4").Value = 108
$ws.Range("V4").Value = 157
$ws.Range("W4").Value = 1.59
$ws.Range("X4").Value = 8.619999999999999
$ws.Range("Y4").Value = 11.27
$ws.Range("Z4").Value = 5.8
$ws.Range("AA4").Value = 112.69
$ws.Range("AB4").Value = 128.58
$ws.Range("AC4").Value = 123
$ws.Range("AD4").Value = 48.25
$ws.Range("AE4").Value = 1209
$ws.Range("AF4").Value = 4.9
$ws.Range("AG4").Value = 21
$ws.Range("AH4").Value = 0.35
$ws.Range("AI4").Value = 16.26
$ws.Range("AJ4").Value = 34458589

# Row 5
$ws.Range("D5").Value = 1120
$ws.Range("E5").Value = 45
$ws.Range("F5").Value = 45
$ws.Range("G5").Value = 51
$ws.Range("H5").Value = 44
$ws.Range("I5").Value = 14
$ws.Range("J5").Value = 30
$ws.Range("K5").Value = 1039
$ws.Range("L5").Value = 629
$ws.Range("M5").Value = 410
$ws.Range("N5").Value = 388
$ws.Range("O5").Value = 21
$ws.Range("P5").Value = 172
$ws.Range("Q5").Value = 98
$ws.Range("R5").Value = -131
$ws.Range("S5").Value = 188
$ws.Range("T5").Value = 69
$ws.Range("U5").Value = 29
$ws.Range("V5").Value = 347
$ws.Range("W5").Value = 4.06
$ws.Range("X5").Value = 3.95
$ws.Range("Y5").Value = 3.63
$ws.Range("Z5").Value = 4.75
$ws.Range("AA5").Value = 153.59
$ws.Range("AB5").Value = 133.03
$ws.Range("AC5").Value = 41
$ws.Range("AD5").Value = 41.34
$ws.Range("AE5").Value = 1186
$ws.Range("AF5").Value = 1.44
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 34458589

# Row 6
$ws.Range("D6").Value = 991
$ws.Range("E6").Value = 46
$ws.Range("F6").Value = 46
$ws.Range("G6").Value = 39
$ws.Range("H6").Value = 15
$ws.Range("I6").Value = -16
$ws.Range("K6").Value = 1062
$ws.Range("L6").Value = 736
$ws.Range("M6").Value = 326
$ws.Range("N6").Value = 312
$ws.Range("P6").Value = 172
$ws.Range("Q6").Value = -112
$ws.Range("R6").Value = 85
$ws.Range("S6").Value = -206
$ws.Range("T6").Value = 6
$ws.Range("U6").Value = -119
$ws.Range("V6").Value = 159
$ws.Range("W6").Value = 4.67
$ws.Range("X6").Value = 1.56
$ws.Range("Y6").Value = -4.5
$ws.Range("Z6").Value = 1.47
$ws.Range("AA6").Value = 225.67
$ws.Range("AB6").Value = 86.90000000000001
$ws.Range("AC6").Value = -46
$ws.Range("AD6").Value = -39.67
$ws.Range("AE6").Value = 953
$ws.Range("AF6").Value = 1.9
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 34458589

# --- Clear now-empty cells in row 6 (columns AG, AH) ---
$ws.Range("AG6:AH6").ClearContents()

# --- Rows 7-9: clear all figures except company name columns A-C ---
$ws.Range("D7:AJ9").ClearContents()
